$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 364, shifting existing rows 364-495 down to 365-496.
$ws.Rows("364:364").Insert()

# Populate the newly inserted row 364 with the new data entry.
$ws.Range("A364").Value = 5
$ws.Range("B364").Value = "Macroferia Regional de Talca"
$ws.Range("C364").Value = "Maule"
$ws.Range("D364").Value = 44988
$ws.Range("E364").Value = 7
$ws.Range("F364").Value = 100112023
$ws.Range("G364").Value = "Brócoli"
$ws.Range("H364").Value = "Sin especificar"
$ws.Range("I364").Value = "Primera"
$ws.Range("J364").Value = 5000
$ws.Range("K364").Value = 700
$ws.Range("L364").Value = 700
$ws.Range("M364").Value = 700
$ws.Range("N364").Value = "$/unidad"
$ws.Range("O364").Value = "Región del Maule"
$ws.Range("P364").Value = 700
$ws.Range("Q364").Value = 1
$ws.Range("R364").Value = "Hortaliza"
